$d = $word.ActiveDocument

function Replace-Text($range, $old, $new) {
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# ---- Title ----
Replace-Text $d.Content "Unraveling Gender Disparity in STEM" "The Fascinating Realm of Chemistry: Unveiling the Secrets of Matter"

# ---- Author name ----
Replace-Text $d.Content "Sarah Miller" "Dennis Wilson"

# ---- Email (paragraph 3) ----
$emailPara = $d.Paragraphs(3).Range
Replace-Text $emailPara "sarahmiller@gmail" "dennis"
Replace-Text $d.Paragraphs(3).Range "com" "wilson@newwave"

# Append the new ".edu" suffix as additional runs after "wilson@newwave"
$findRng = $d.Content
$findRng.Find.Execute("wilson@newwave", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $findRng.End
$insertRng = $d.Range($endPos, $endPos)
$insertRng.InsertAfter(".")
$insertRng2 = $d.Range($endPos + 1, $endPos + 1)
$insertRng2.InsertAfter("edu")

# ---- Body paragraph (paragraph 5, the main essay paragraph) ----
$bodyPara = $d.Paragraphs(5).Range

Replace-Text $bodyPara `
    "Throughout history, the realm of science, technology, engineering, and mathematics (STEM) has been predominantly male-dominated" `
    "In the vast tapestry of scientific disciplines, Chemistry stands as a beacon of discovery, shedding light on the intricate workings of matter"

Replace-Text $bodyPara `
    " This disparity is a global phenomenon, transcending cultural and socioeconomic boundaries" `
    " It is a subject that touches every aspect of our lives, from the food we consume to the air we breathe, and even the clothes we wear"

Replace-Text $bodyPara `
    " Consequently, society has missed out on the invaluable contributions of women in STEM fields, resulting in a skewed representation of perspectives, ideas, and innovations. This essay delves into the multifaceted issue of gender disparity in STEM, examining its root causes, detrimental effects, and potential solutions to foster a more inclusive environment." `
    " In this exploration, we will delve into the captivating world of Chemistry, revealing the fundamental concepts that govern the behavior of substances and the remarkable transformations they undergo"

Replace-Text $bodyPara `
    "Gender stereotypes and societal expectations play a pivotal role in shaping girls' and women's career choices" `
    "Within the realm of Chemistry, we will unravel the enigmatic nature of atoms and molecules, the building blocks of all matter"

Replace-Text $bodyPara `
    " From an early age, girls are often discouraged from pursuing careers in STEM, as these fields are traditionally perceived as masculine" `
    " We will explore the periodic table, a roadmap that unveils the properties and relationships of these fundamental particles"

Replace-Text $bodyPara `
    " This perception is reinforced by the limited visibility of female role models in STEM, perpetuating the notion that these fields are not suitable for women. Furthermore, the lack of gender diversity in STEM creates a hostile environment for women, where they face discrimination, prejudice, and a lack of support." `
    " Through engaging experiments and captivating demonstrations, we will witness the symphony of chemical reactions, marveling at the vibrant colors and intriguing changes that accompany them"

Replace-Text $bodyPara `
    "The consequences of gender disparity in STEM are multi-faceted and far-reaching" `
    "Furthermore, we will probe the depths of chemical bonding, the forces that hold atoms together and determine the properties of substances"

Replace-Text $bodyPara `
    " It deprives society of the talents and contributions of a large pool of potential scientists, engineers, and innovators" `
    " We will delve into the mysteries of acids and bases, unveiling their roles in everyday phenomena and their applications in various industries"

Replace-Text $bodyPara `
    " This has implications for economic growth, as well as the development of new technologies and solutions to address global challenges. Moreover, the underrepresentation of women in STEM reinforces gender stereotypes and perpetuates the cycle of discrimination, creating a vicious loop that is difficult to break" `
    " As we progress in our journey through Chemistry, we will appreciate the intricate dance of particles and the elegance of chemical principles that orchestrate the material world around us"

# ---- Summary paragraph (paragraph 7) ----
$summaryPara = $d.Paragraphs(7).Range

Replace-Text $summaryPara `
    "The gender disparity in STEM is a multifaceted issue with profound implications for society" `
    "Chemistry, a captivating realm of scientific inquiry, uncovers the secrets of matter and its transformations"

Replace-Text $summaryPara `
    " Rooted in societal stereotypes and a lack of female role models, it results in a hostile environment for women in STEM, leading to discrimination, prejudice, and a lack of support" `
    " Through the study of atoms, molecules, and chemical reactions, we gain profound insights into the behavior of substances and their applications in various industries"

Replace-Text $summaryPara `
    " This disparity has detrimental consequences, including the loss of talent, the perpetuation of gender stereotypes, and the stifling of innovation. To address this issue, comprehensive efforts are required to challenge stereotypes, promote female role models, create inclusive environments, and implement policies that support women in STEM. By fostering a more diverse and inclusive STEM workforce, society can unlock the full potential of innovation and progress" `
    " By delving into the fundamental principles of Chemistry, we not only enhance our understanding of the world around us but also equip ourselves with valuable tools for solving real-world problems, paving the way for future scientific advancements and technological breakthroughs"

# ---- Add a new empty paragraph at the very end of the document body ----
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
